# "Updates to Team Management and detection setup"
#
# 1. Sheet "Before covid": the old meeting-notes block (rows 34:44) is
#    removed, leaving a single blank (but still styled) row 34.
# 2. Sheet "After Covid": the task rows are relabelled/reorganised and two
#    new notes ("bulk of writing" / "have diagram prepared") are added,
#    plus three new columns (F:H) get explicit widths.
# 3. Selections on both sheets move to reflect where the author was last
#    working.

$wb = $excel.ActiveWorkbook
$wsBefore = $wb.Worksheets.Item("Before covid")
$wsAfter  = $wb.Worksheets.Item("After Covid")

# ------------------------------------------------------------------
# "Before covid" sheet: wipe out the old meeting-log rows 34-44, leaving
# row 34 as a blank (but still styled) placeholder row.
# ------------------------------------------------------------------
$wsBefore.Range("A34:L44").ClearContents()

# ------------------------------------------------------------------
# "After Covid" sheet: new column widths for F, G, H
# ------------------------------------------------------------------
$wsAfter.Columns("F").ColumnWidth = 16.0
$wsAfter.Columns("G").ColumnWidth = 20.666666666666668
$wsAfter.Columns("H").ColumnWidth = 13.666666666666666

# ------------------------------------------------------------------
# "After Covid" sheet: relabel column A task names
# ------------------------------------------------------------------
$wsAfter.Range("A2").Value  = "date"
$wsAfter.Range("A3").Value  = "Ideation"
$wsAfter.Range("A4").Value  = "Solution"
$wsAfter.Range("A5").Value  = "Baseline"
$wsAfter.Range("A6").Value  = "Testing"
$wsAfter.Range("A7").Value  = "Management"
$wsAfter.Range("A8").Value  = "Mathematical"
$wsAfter.Range("A9").Value  = "API"
$wsAfter.Range("A10").Value = "GUI"
$wsAfter.Range("A11").Value = "Iterative"

# Column I note, same for every data row
$wsAfter.Range("I3").Value  = "final touches"
$wsAfter.Range("I4").Value  = "final touches"
$wsAfter.Range("I5").Value  = "final touches"
$wsAfter.Range("I6").Value  = "final touches"
$wsAfter.Range("I7").Value  = "final touches"
$wsAfter.Range("I8").Value  = "final touches"
$wsAfter.Range("I9").Value  = "final touches"
$wsAfter.Range("I10").Value = "final touches"
$wsAfter.Range("I11").Value = "final touches"

# Row 3: Ideation notes
$wsAfter.Range("C3").Value = "bulk of writing"

# Row 4: Solution notes
$wsAfter.Range("G4").Value = "have diagram prepared"
$wsAfter.Range("H4").Value = "bulk of writing"

# Row 5: Baseline notes
$wsAfter.Range("B5").Value = "find new dataset"
$wsAfter.Range("C5").Value = "setup detection in one file"
$wsAfter.Range("D5").Value = "setup environment"
$wsAfter.Range("E5").Value = "setup training"

# Row 8: Mathematical notes
$wsAfter.Range("C8").Value = "evaluate detection"
$wsAfter.Range("E8").Value = "evaluate training"

# ------------------------------------------------------------------
# Selections: "Before covid" ends up with H42 selected (not activated,
# so it stays off the tab-selected sheet); "After Covid" keeps the
# tab-selected flag and ends with G8 selected.
# ------------------------------------------------------------------
$wsBefore.Range("H42").Select()
$wsAfter.Range("G8").Select()
